# Update cryptos list values per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'64.269.31"
$ws.Range('E2').Value = "'  +0.54%  "
$ws.Range('D3').Value = "'3.325.45"
$ws.Range('E3').Value = "'  -0.01%  "
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = "'  -0.06%  "
$ws.Range('D5').Value = "'552.99"
$ws.Range('E5').Value = "'  +0.42%  "
$ws.Range('D6').Value = "'173.25"
$ws.Range('E6').Value = "'  +0.76%  "
$ws.Range('D7').Value = "'0.619"
$ws.Range('E7').Value = "'  +1.09%  "
$ws.Range('D8').Value = "'0.999"
$ws.Range('E8').Value = "'  -0.01%  "
$ws.Range('D9').Value = "'3.317.58"
$ws.Range('E9').Value = "'  -0.18%  "
$ws.Range('D10').Value = "'0.170"
$ws.Range('E10').Value = "'  +5.85%  "
$ws.Range('D11').Value = "'0.630"
$ws.Range('E11').Value = "'  +1.59%  "
$ws.Range('D12').Value = "'53.33"
$ws.Range('E12').Value = "'  +0.66%  "
$ws.Range('E13').Value = "'  +3.14%  "
$ws.Range('D14').Value = "'9.06"
$ws.Range('E14').Value = "'  +0.88%  "
$ws.Range('D15').Value = "'3.849.45"
$ws.Range('E15').Value = "'  -0.20%  "
$ws.Range('E16').Value = "'  +2.96%  "
$ws.Range('D17').Value = "'18.09"
$ws.Range('E17').Value = "'  -0.71%  "
$ws.Range('D18').Value = "'3.317.62"
$ws.Range('E18').Value = "'  +0.10%  "
$ws.Range('D19').Value = "'64.317.68"
$ws.Range('E19').Value = "'  +0.79%  "
$ws.Range('D20').Value = "'11.69"
$ws.Range('E20').Value = "'  -0.33%  "
$ws.Range('D21').Value = "'0.984"
$ws.Range('E21').Value = "'  +1.55%  "
$ws.Range('D22').Value = "'453.86"
$ws.Range('E22').Value = "'  +6.24%  "
$ws.Range('D23').Value = "'5.12"
$ws.Range('E23').Value = "'  +9.38%  "
$ws.Range('D24').Value = "'4.06"
$ws.Range('E24').Value = "'  -0.24%  "
$ws.Range('B25').Value = "'Litecoin"
$ws.Range('C25').Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range('D25').Value = "'87.13"
$ws.Range('E25').Value = "'  +3.74%  "
$ws.Range('B26').Value = "'InternetComputer(DFINITY)"
$ws.Range('C26').Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range('D26').Value = "'13.92"
$ws.Range('E26').Value = "'  +4.46%  "
$ws.Range('D27').Value = "'2.86"
$ws.Range('E27').Value = "'  +1.77%  "
$ws.Range('D28').Value = "'10.61"
$ws.Range('E28').Value = "'  +0.15%  "
$ws.Range('D29').Value = "'30.90"
$ws.Range('E29').Value = "'  +4.46%  "
$ws.Range('D30').Value = "'8.56"
$ws.Range('E30').Value = "'  +0.16%  "
$ws.Range('D31').Value = "'6.51"
$ws.Range('E31').Value = "'  -2.04%  "
$ws.Range('D32').Value = "'11.39"
$ws.Range('E32').Value = "'  +0.22%  "
$ws.Range('D33').Value = "'61.54"
$ws.Range('E33').Value = "'  +5.97%  "
$ws.Range('D34').Value = "'564.76"
$ws.Range('E34').Value = "'  -4.87%  "
$ws.Range('D35').Value = "'0.107"
$ws.Range('E35').Value = "'  +0.13%  "
$ws.Range('E36').Value = "'  +0.06%  "
$ws.Range('B37').Value = "'Kaspa"
$ws.Range('C37').Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range('D37').Value = "'0.140"
$ws.Range('E37').Value = "'  -1.43%  "
$ws.Range('B38').Value = "'Stacks"
$ws.Range('C38').Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range('D38').Value = "'3.51"
$ws.Range('E38').Value = "'  +1.88%  "
$ws.Range('D39').Value = "'35.21"
$ws.Range('E39').Value = "'  +0.04%  "
$ws.Range('D40').Value = "'0.365"
$ws.Range('E40').Value = "'  +0.53%  "
$ws.Range('D41').Value = "'0.0₃0727"
$ws.Range('E41').Value = "'  -2.29%  "
$ws.Range('D42').Value = "'3.050.77"
$ws.Range('E42').Value = "'  -1.11%  "
$ws.Range('D43').Value = "'0.0416"
$ws.Range('E43').Value = "'  +2.84%  "
$ws.Range('E44').Value = "'  -1.37%  "
$ws.Range('D45').Value = "'3.20"
$ws.Range('E45').Value = "'  +0.70%  "
$ws.Range('E46').Value = "'  +0.87%  "
$ws.Range('E47').Value = "'  +3.36%  "
$ws.Range('D48').Value = "'0.999"
$ws.Range('E48').Value = "'  -0.04%  "
$ws.Range('D49').Value = "'140.53"
$ws.Range('E49').Value = "'  +6.23%  "
$ws.Range('D50').Value = "'2.50"
$ws.Range('E50').Value = "'  -3.12%  "
$ws.Range('D51').Value = "'8.13"
$ws.Range('E51').Value = "'  +0.04%  "
